# code-complete-status.xlsx edit:
#  - add links to chapters 23.1-23.3 to lab 06 (column E)
#  - add link to chapter 19.4 to lab 05 (column D)
#
# Concretely:
#  * row "19.4" (old row 30) gets a checkmark in column D (lab 05)
#  * row "23.1" (old row 31) gets a checkmark in column E (lab 06)
#  * row "23.2" (old row 32) gets a checkmark in column E (lab 06)
#  * a brand new row "23.3" is inserted after "23.2" with a checkmark in
#    column E (lab 06) as well
#  * the two rows that previously had no "used?" (N) formula (16.1 / row 25,
#    19.3 / row 29) get the same shared formula as the rest of the column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fill the two existing rows that were missing the N formula ---------
$ws.Range("N25").Formula = "=IF(SUM(B25:M25)>0,1,0)"
$ws.Range("N29").Formula = "=IF(SUM(B29:M29)>0,1,0)"

# --- mark chapter 19.4 (row 30) as linked from lab 05 (column D) --------
$ws.Range("D30").Value = 1

# --- mark chapters 23.1 and 23.2 (rows 31-32) as linked from lab 06 -----
$ws.Range("E31").Value = 1
$ws.Range("E32").Value = 1

# --- insert a new row for chapter 23.3 right after 23.2 (row 33) --------
$ws.Rows(33).Insert()
$ws.Range("A33").Value = "23.3"
$ws.Range("E33").Value = 1
$ws.Range("N33").Formula = "=IF(SUM(B33:M33)>0,1,0)"

# --- extend the conditional formatting range to include the new row -----
$cf = $ws.Range("N3:N37").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("N3:N38"))

# --- move the selection like in the edited workbook ----------------------
$ws.Range("I29").Select()
